# Continue the Working Time Log: a new shift was clocked for the same day as
# row 9 (43054 = the date serial) - entered in row 10 (A10:C10). The "Work
# Time" column (D) is a shared formula already present for the whole table,
# so it recalculates on its own once B10/C10 are filled in; the table total
# in D36 recalculates along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number formats / styles) from the row above so the
# new row matches the existing "Date" / "Start Time" / "End Time" look
# exactly, instead of mutating each cell's NumberFormat by hand (which would
# mint brand-new, duplicate style entries rather than reusing the ones the
# sheet already has).
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C9").Copy()
$ws.Range("C10").PasteSpecial(-4122)   # xlPasteFormats

# New clock-in / clock-out entry.
$ws.Range("A10").Value = 43054
$ws.Range("B10").Value = 0.45833333333333331
$ws.Range("C10").Value = 0.4861111111111111

# Leave the cursor where the user moved it next.
[void]$ws.Range("C11").Select()
